$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.230.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.484.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.556"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.519.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.935.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.151.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.504.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "323.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.12%  "
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.984"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0750"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  +3.03%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.993"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.780"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "277.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("E43").Value = "  +2.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.598"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0922"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0500"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  +0.14%  "
